# Refresh the cryptos price/volume columns (and re-rank dogwifhat vs OKB)
# to match the latest scrape, per the automated GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.612.61"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "2.604.63"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "571.97"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").Value = "155.14"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").Value = "0.118"
$ws.Range("E9").Value = "  -6.96%  "
$ws.Range("D10").Value = "5.83"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "0.382"
$ws.Range("E11").Value = "  -5.05%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "28.25"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").Value = "3.083.03"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  -8.20%  "
$ws.Range("D16").Value = "63.413.00"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").Value = "2.591.83"
$ws.Range("E17").Value = "  -4.16%  "
$ws.Range("D18").Value = "11.99"
$ws.Range("E18").Value = "  -5.32%  "
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "4.54"
$ws.Range("E20").Value = "  -5.49%  "
$ws.Range("D21").Value = "342.83"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "67.22"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "1.79"
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").Value = "0.0000108"
$ws.Range("E25").Value = "  -3.82%  "
$ws.Range("D26").Value = "591.54"
$ws.Range("E26").Value = "  +4.35%  "
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  -5.23%  "
$ws.Range("E28").Value = "  -4.21%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").Value = "1.73"
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("D34").Value = "6.53"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("D35").Value = "5.43"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").Value = "0.404"
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "19.69"
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D39").Value = "154.93"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("E41").Value = "  +0.02%  "

# Rows 42/43: dogwifhat overtook OKB in rank, so the two rows swap coin identity.
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  +7.39%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.30"
$ws.Range("E43").Value = "  -3.38%  "

$ws.Range("D44").Value = "156.34"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").Value = "3.91"
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("D46").Value = "23.11"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "0.0587"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").Value = "0.629"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").Value = "0.0247"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "18.78"
$ws.Range("E51").Value = "  -5.60%  "
